# Change all spreadsheet timeout values to 6 or 7 seconds vice 6000 or 7000
# seconds, and add a comment to spreadsheets where timeout values are
# specified that timeout values are seconds, not milliseconds.
# Also fix "Plateform"/"plateform" typos to "Platform"/"platform" in the
# config-comment column.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "c-demo_ui"
$ws2 = $wb.Worksheets.Item(2)   # "t-excel_excel"

# --- waitTime: 7000 (ms) -> 7 (seconds); add explanatory comment ---
$ws1.Range("B23").Value = "7"

# --- Fix "Plateform type for remote web driver intializing" typo ---
$ws1.Range("C7").Value  = "Platform type for remote web driver intializing"
$ws1.Range("C9").Value  = "Platform type for remote web driver intializing"
$ws1.Range("C11").Value = "Platform type for remote web driver intializing"
$ws1.Range("C13").Value = "Platform type for remote web driver intializing"

# --- Fix "Version for plateform type selected" typo ---
$ws1.Range("C8").Value  = "Version for platform type selected"
$ws1.Range("C10").Value = "Version for platform type selected"
$ws1.Range("C12").Value = "Version for platform type selected"
$ws1.Range("C14").Value = "Version for platform type selected"

$ws1.Range("C23").Value = "Wait time is in seconds, not milliseconds"

# --- View state: make "c-demo_ui" the active/selected sheet & cell ---
$ws1.Activate()
$ws1.Range("C29").Select()
